# Auto-generated edit script: updates cryptos list price/volume figures
# per the commit 'Updated cryptos list ... with GitHub Actions'.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellRef,
        [string]$NewValue
    )
    $range = $ws.Range($CellRef)
    # Force text interpretation so numeric-looking strings (e.g. "571.65",
    # "1.00", "0.999") are NOT silently coerced into numbers - these
    # source cells are plain text in the workbook.
    $range.NumberFormat = "@"
    $range.Value = $NewValue
    # Restore the default ("Normal") cell style so we don't leave a
    # lingering custom number-format style behind.
    $range.Style = "Normal"
}

Set-TextValue 'D2' '60.959.43'
Set-TextValue 'E2' '  +0.18%  '
Set-TextValue 'D3' '3.385.56'
Set-TextValue 'E3' '  -0.53%  '
Set-TextValue 'E4' '  +0.00%  '
Set-TextValue 'D5' '571.65'
Set-TextValue 'E5' '  -0.25%  '
Set-TextValue 'D6' '142.19'
Set-TextValue 'E6' '  -0.11%  '
Set-TextValue 'E7' '  -0.03%  '
Set-TextValue 'E8' '  -0.06%  '
Set-TextValue 'E9' '  +0.44%  '
Set-TextValue 'E10' '  -1.91%  '
Set-TextValue 'E11' '  -1.40%  '
Set-TextValue 'D12' '3.964.90'
Set-TextValue 'E12' '  -0.49%  '
Set-TextValue 'E13' '  +1.73%  '
Set-TextValue 'D14' '27.76'
Set-TextValue 'E14' '  -2.11%  '
Set-TextValue 'E15' '  -0.82%  '
Set-TextValue 'D16' '3.369.41'
Set-TextValue 'E16' '  -1.16%  '
Set-TextValue 'D17' '61.062.42'
Set-TextValue 'E17' '  +0.15%  '
Set-TextValue 'E18' '  -3.69%  '
Set-TextValue 'D19' '13.67'
Set-TextValue 'E19' '  -4.07%  '
Set-TextValue 'E20' '  -2.34%  '
Set-TextValue 'D21' '383.59'
Set-TextValue 'E21' '  -1.63%  '
Set-TextValue 'D22' '74.95'
Set-TextValue 'E22' '  +2.70%  '
Set-TextValue 'E23' '  -2.62%  '
Set-TextValue 'E24' '  +0.35%  '
Set-TextValue 'E25' '  -4.97%  '
Set-TextValue 'D26' '3.525.27'
Set-TextValue 'E26' '  -0.63%  '
Set-TextValue 'D27' '0.182'
Set-TextValue 'E27' '  +0.99%  '
Set-TextValue 'D28' '0.999'
Set-TextValue 'E28' '  -0.08%  '
Set-TextValue 'D29' '7.29'
Set-TextValue 'E29' '  -1.59%  '
Set-TextValue 'E30' '  -2.19%  '
Set-TextValue 'E31' '  -0.66%  '
Set-TextValue 'B32' 'USDe'
Set-TextValue 'C32' 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue 'D32' '1.00'
Set-TextValue 'E32' '  -0.02%  '
Set-TextValue 'B33' 'Fetch.AI'
Set-TextValue 'C33' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 'D33' '1.40'
Set-TextValue 'E33' '  -4.10%  '
Set-TextValue 'D34' '23.26'
Set-TextValue 'E34' '  -2.56%  '
Set-TextValue 'E35' '  -0.95%  '
Set-TextValue 'D36' '166.98'
Set-TextValue 'E36' '  -0.17%  '
Set-TextValue 'D37' '3.418.40'
Set-TextValue 'E37' '  -0.39%  '
Set-TextValue 'D38' '5.00'
Set-TextValue 'E38' '  -2.28%  '
Set-TextValue 'D39' '1.47'
Set-TextValue 'E39' '  -4.72%  '
Set-TextValue 'E40' '  -2.35%  '
Set-TextValue 'D41' '26.64'
Set-TextValue 'E41' '  -1.93%  '
Set-TextValue 'E42' '  -0.02%  '
Set-TextValue 'D43' '0.779'
Set-TextValue 'E43' '  -1.39%  '
Set-TextValue 'E44' '  -2.49%  '
Set-TextValue 'E45' '  -1.98%  '
Set-TextValue 'D46' '1.14'
Set-TextValue 'E46' '  -0.03%  '
Set-TextValue 'D47' '2.460.74'
Set-TextValue 'E47' '  -4.48%  '
Set-TextValue 'D48' '23.06'
Set-TextValue 'E48' '  -0.33%  '
Set-TextValue 'E49' '  -2.83%  '
Set-TextValue 'D50' '2.16'
Set-TextValue 'E50' '  +7.04%  '
Set-TextValue 'E51' '  +1.14%  '
